$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)
$ws.Range("A1").Value = "Test"
Write-Output "done"
